$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Benchmark results updates ---

# Remove the stale #DIV/0! helper formula in C9 (division by the
# not-yet-populated 6MHz hardware-frequency column).
$ws.Range("C9").ClearContents()

# Add a new label for the 12MHz average-CPI row, using a new shared string.
$ws.Range("A11").Value = "average CPI on hardware at 12MHz"

# The 6MHz-era scratch rows (13-16) are no longer the focus of the sheet,
# so hide them, while keeping the 12MHz results (rows above) visible.
$ws.Rows.Item(13).Hidden = $true
$ws.Rows.Item(14).Hidden = $true
$ws.Rows.Item(15).Hidden = $true
$ws.Rows.Item(16).Hidden = $true

# Column A needs to be much wider to fit the longer labels; column B holds
# supporting data that's no longer displayed, so hide it.
$ws.Columns.Item(1).ColumnWidth = 42
$ws.Columns.Item(2).Hidden = $true

# Move/update the active selection to where the author left off reviewing.
$ws.Range("A24").Select() | Out-Null
